# ---------------------------------------------------------------------------
# Applies the "2nd VM Exercise.docx" edit described by the commit diff:
#   * retitle the Learning Objective line
#   * rewrite the Context paragraph and drop the blank line that followed it
#   * turn "Instructions:" into a proper Heading-4 "Steps to be executed by
#     the student:" lead-in (with shading/spacing + a bookmark), restoring
#     the direct character formatting that changing the paragraph style
#     resets
#   * a handful of in-place wording tweaks further down the list
#   * add a left/first-line indent to the blank paragraph that trails the
#     last bullet
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# 1) Title line -------------------------------------------------------------
$d.Content.Find.Execute(
    " Cleaning Data ", $false, $false, $false, $false, $false, $true, 1, $false,
    " Cleaning and preparing data for analysis", 2)

# 2) Context paragraph text --------------------------------------------------
$oldContext = " Real world data usually comes from various sources. If measures are not put in place to authenticate the data entered and stored, it may come with unwanted characters, improper fonts and cases, nulls, duplicates  and so on. It is important that an analyst cleans data before proceeding to the analysis stage. In this exercise, the learner is going to remove duplicates and clean data."
$newContext = " In this exercise, the student is going to identify/remove duplicates using the Unique Tool and modify the case of values in a column using the Data Cleansing Tool. Real world data usually comes in dirty. If measures are not put in place to authenticate the data as it is entered and stored, it may come with unwanted characters, improper fonts, missing values, and duplicates. It is important that an analyst cleans the data before proceeding to the analysis stage."
$d.Content.Find.Execute($oldContext, $false, $false, $false, $false, $false, $true, 1, $false, $newContext, 2)

# 3) Drop the blank paragraph that used to sit between "Context" and
#    "Instructions" (paragraph #4 at this point; nothing above touched the
#    paragraph count).
$d.Paragraphs.Item(4).Range.Delete()

# 4) "Instructions:" -> "Steps to be executed by the student:" heading -----
$instr = $d.Paragraphs.Item(4)

# 4a) swap the text while the run still carries its original direct
#     formatting (rFonts/bold/sz/szCs/rtl) so Find can match cleanly.
$instr.Range.Find.Execute(
    "Instructions: ", $false, $false, $false, $false, $false, $true, 1, $false,
    "Steps to be executed by the student: ", 2)

# 4b) promote the paragraph to the Heading 4 style (this resets direct
#     character formatting on its runs + paragraph mark, restored below).
$instr.Style = "Heading 4"

# 4c) paragraph-level formatting added alongside the style change.
$instr.Format.KeepWithNext = $false
$instr.Format.KeepTogether = $false
$instr.Format.Shading.Texture = 0
$instr.Format.Shading.BackgroundPatternColor = 16777215
$instr.Format.SpaceBefore = 18
$instr.Format.SpaceAfter = 12
$instr.Format.LineSpacingRule = 5
$instr.Format.LineSpacing = 15

# 4d) restore the paragraph-mark run formatting (rFonts/bold/sz/szCs) that
#     the style switch wiped - matches the unchanged <w:pPr><w:rPr> block in
#     the target XML.
$markRng = $instr.Range.Duplicate()
$markRng.Collapse(0)
$markRng.Font.NameAscii = "Poppins"
$markRng.Font.Name = "Poppins"
$markRng.Font.NameFarEast = "Poppins"
$markRng.Font.NameBi = "Poppins"
$markRng.Font.Bold = $true
$markRng.Font.Size = 10

# 4e) restore + extend the run formatting on "Steps to be executed by the
#     student: " (rFonts/bold/color/sz/szCs/rtl).
$runRng = $instr.Range.Duplicate()
$runRng.Font.NameAscii = "Poppins"
$runRng.Font.Name = "Poppins"
$runRng.Font.NameFarEast = "Poppins"
$runRng.Font.NameBi = "Poppins"
$runRng.Font.Bold = $true
$runRng.Font.Size = 10
$runRng.Font.Color = 0

# 4f) empty bookmark at the very start of the paragraph.
$bmRng = $instr.Range.Duplicate()
$bmRng.Collapse(1)
$d.Bookmarks.Add("_ii7v72xogp40", $bmRng)

# 5) Wording tweaks further down the list -----------------------------------
$oldStep8 = "Find and drag the Data Cleaning tool and connect it to the " + [string][char]0x201C + "U" + [string][char]0x201D + " output anchor of the Unique tool"
$newStep8 = "Find and drag the Data Cleansing tool and connect it to the " + [string][char]0x201C + "U" + [string][char]0x201D + " output anchor of the Unique tool from the step above."
$d.Content.Find.Execute($oldStep8, $false, $false, $false, $false, $false, $true, 1, $false, $newStep8, 2)

$d.Content.Find.Execute("Leave all other checkboxes", $false, $false, $false, $false, $false, $true, 1, $false, "Leave other checkboxes", 2)

$oldToggle = "Toggle between the input and output anchors to check in Country/Region in is the title case"
$newToggle = "Toggle between the input and output anchors to check in Country/Region changed from lower case to title case"
$d.Content.Find.Execute($oldToggle, $false, $false, $false, $false, $false, $true, 1, $false, $newToggle, 2)

# 6) Indent the blank paragraph right after the "Toggle..." bullet ----------
$afterToggle = $d.Paragraphs.Item(11)
$afterToggle.Format.LeftIndent = 36
$afterToggle.Format.FirstLineIndent = 0
